$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells that contain list-like placeholder strings for duplicated teacher
# entries should be collapsed to a simple "-" (removing the duplication).
$cells = @("E3", "E4", "E6", "E7", "D11", "D12", "D14", "D15", "B18", "B19", "B20", "B21", "C18", "C19", "C20", "C21")

foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = "-"
}
